# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404
#   *_new -> *_FV2410
# Then wrap the data range in a table (Table1) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $colOld = $i + 1
    $colNew = $i + 12
    $ws.Cells.Item(1, $colOld).Value = $baseNames[$i] + "_FV2404"
    $ws.Cells.Item(1, $colNew).Value = $baseNames[$i] + "_FV2410"
}

# Turn A1:U57 into an Excel table (ListObject) that uses the header row values
# as the column names, with an auto filter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1, top-left of the scrolling pane is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
